# Update "想去人数" (F column) counts on three sheets to match the
# gh-pages regeneration snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 392
$ws1.Range("F9").Value  = 205
$ws1.Range("F12").Value = 1060
$ws1.Range("F15").Value = 197
$ws1.Range("F16").Value = 1528
$ws1.Range("F22").Value = 1164
$ws1.Range("F24").Value = 1913
$ws1.Range("F25").Value = 2681
$ws1.Range("F26").Value = 1468
$ws1.Range("F27").Value = 69
$ws1.Range("F28").Value = 49
$ws1.Range("F29").Value = 445
$ws1.Range("F30").Value = 599
$ws1.Range("F31").Value = 1319
$ws1.Range("F32").Value = 835
$ws1.Range("F33").Value = 1410
$ws1.Range("F36").Value = 795
$ws1.Range("F37").Value = 647
$ws1.Range("F38").Value = 691
$ws1.Range("F39").Value = 873
$ws1.Range("F40").Value = 374
$ws1.Range("F41").Value = 260

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 648
$ws2.Range("F23").Value = 16

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 392
$ws4.Range("F15").Value = 205
$ws4.Range("F19").Value = 1060
$ws4.Range("F21").Value = 197
$ws4.Range("F22").Value = 1528
$ws4.Range("F29").Value = 1164
$ws4.Range("F30").Value = 2681
$ws4.Range("F31").Value = 1468
$ws4.Range("F32").Value = 69
$ws4.Range("F35").Value = 445
$ws4.Range("F36").Value = 1319
$ws4.Range("F39").Value = 835
$ws4.Range("F40").Value = 1410
$ws4.Range("F41").Value = 795
$ws4.Range("F42").Value = 647
$ws4.Range("F43").Value = 691
$ws4.Range("F44").Value = 873
$ws4.Range("F45").Value = 374
$ws4.Range("F47").Value = 16
$ws4.Range("F48").Value = 260
